$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.330.31'
$ws.Range('E2').Value = '  +4.93%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.587.11'
$ws.Range('E3').Value = '  +6.31%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '505.90'
$ws.Range('E5').Value = '  +2.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.36'
$ws.Range('E6').Value = '  -1.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.579'
$ws.Range('E8').Value = '  -4.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.619.68'
$ws.Range('E9').Value = '  +6.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.59'
$ws.Range('E10').Value = '  +5.10%  '
$ws.Range('E11').Value = '  +2.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.343'
$ws.Range('E12').Value = '  +2.57%  '
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.037.49'
$ws.Range('E14').Value = '  +6.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.330.92'
$ws.Range('E15').Value = '  +4.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.72'
$ws.Range('E16').Value = '  +3.80%  '
$ws.Range('E17').Value = '  +3.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.613.94'
$ws.Range('E18').Value = '  +7.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.80'
$ws.Range('E19').Value = '  +2.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '345.74'
$ws.Range('E20').Value = '  +6.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.38'
$ws.Range('E21').Value = '  +2.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.09'
$ws.Range('E22').Value = '  +2.55%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.29'
$ws.Range('E24').Value = '  +3.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.422'
$ws.Range('E25').Value = '  +4.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.165'
$ws.Range('E26').Value = '  +3.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.700.03'
$ws.Range('E27').Value = '  +7.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.992'
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0857'
$ws.Range('E29').Value = '  +6.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.46'
$ws.Range('E30').Value = '  +1.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '156.05'
$ws.Range('E32').Value = '  +3.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.28'
$ws.Range('E33').Value = '  +2.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.56'
$ws.Range('E34').Value = '  +1.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.75'
$ws.Range('E35').Value = '  +7.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.01'
$ws.Range('E36').Value = '  +6.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.20'
$ws.Range('E37').Value = '  +4.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.858'
$ws.Range('E38').Value = '  +24.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.851'
$ws.Range('E39').Value = '  +2.96%  '
$ws.Range('E40').Value = '  +5.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.77'
$ws.Range('E41').Value = '  +5.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '299.30'
$ws.Range('E42').Value = '  +7.22%  '
$ws.Range('E43').Value = '  +3.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0569'
$ws.Range('E44').Value = '  +5.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.624'
$ws.Range('E45').Value = '  +3.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0999'
$ws.Range('E46').Value = '  -0.63%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.994'
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.89'
$ws.Range('E48').Value = '  +10.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.97'
$ws.Range('E49').Value = '  +6.69%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.050.85'
$ws.Range('E50').Value = '  +8.28%  '
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0235'
$ws.Range('E51').Value = '  +1.73%  '
